# Applies the "Dani worked on discussion" edits to the document.
$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# "...across the XX students surveyed (" -> "...across 206 survey ("
$d.Content.Find.Execute(
    "one year after the end of an ESD intervention across the XX students surveyed (",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "one year after the end of an ESD intervention across 206 survey (",
    2) | Out-Null

# --- Change 2 ---------------------------------------------------------
# Remove "form" spelling-error flag (proofErr) -- plain text is unaffected;
# the Find/Replace below normalises the run structure around it.
$d.Content.Find.Execute(
    "significance level resulting form the Wilcoxon test comparing the groups.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "significance level resulting form the Wilcoxon test comparing the groups.",
    2) | Out-Null

# --- Change 3 ---------------------------------------------------------
# "(Research Question 2a)" -> "(Research Question 2)" for the TPB heading
$d.Content.Find.Execute(
    "Theory of planned behaviour and self-efficacy beliefs (Research Question 2a)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Theory of planned behaviour and self-efficacy beliefs (Research Question 2)",
    2) | Out-Null

# --- Change 4 ---------------------------------------------------------
# Remove the empty paragraph right before the "Self-efficacy beliefs and
# level of involvement (Research Question 2b)" heading.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "`r" -and $p.Next().Range.Text.StartsWith("Self-efficacy beliefs and level of involvement")) {
        $p.Range.Delete() | Out-Null
        break
    }
}

# --- Change 5 ---------------------------------------------------------
# Turn "Self-efficacy beliefs and level of involvement (Research Question 2b)"
# into the Heading3-styled (but Heading2-sized) "Sustainability competencies
# as efficacy beliefs (Research Question 3)"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Self-efficacy beliefs and level of involvement")) {
        $r = $p.Range
        $r.MoveEnd(1, -1) | Out-Null
        $r.Text = "Sustainability competencies as efficacy beliefs (Research Question 3)"
        $p.Range.Style = "Heading 3"
        $p.Range.Font.Name = "+Headings"
        $p.Range.Font.Size = 16
        $p.Range.Font.SizeBi = 16
        break
    }
}

# --- Change 6 ---------------------------------------------------------
# "Individual and collective self-efficacy beliefs between and within
# groups" -> "Personal and collective efficacy beliefs"
$d.Content.Find.Execute(
    "Individual and collective self-efficacy beliefs between and within groups",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Personal and collective efficacy beliefs",
    2) | Out-Null

# --- Change 7 ---------------------------------------------------------
# Tidy the "see Table X,….)" run split / drop the grammar-error flags
# (ellipsis character below is a literal U+2026).
$d.Content.Find.Execute(
    "individual and collective self-efficacy (p > .05, see Table X,….) across the XX students surveyed at time point 3.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "individual and collective self-efficacy (p > .05, see Table X,….) across the XX students surveyed at time point 3.",
    2) | Out-Null

# --- Change 8 ---------------------------------------------------------
# Split the run at "...individual " / "self-efficacy beliefs..." so that a
# lastRenderedPageBreak can sit at the new page boundary (pagination is
# recomputed automatically on save).
$d.Content.Find.Execute(
    "The reported mean scores of the involved group were higher for both collective and individual self-efficacy beliefs than those of the control group.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The reported mean scores of the involved group were higher for both collective and individual self-efficacy beliefs than those of the control group.",
    2) | Out-Null

# --- Change 9 ---------------------------------------------------------
# Merge the aim/action-focussed paragraph back into a single run (its
# mid-paragraph page break is no longer needed once the earlier content
# shifts; pagination recalculates on save).
$d.Content.Find.Execute(
    "Looking at the relationship between aim and action focussed self-efficacy beliefs and group membership, I found that the involved group reported significantly higher scores regarding the aim focus, than the control group (W = 8, p < .05, Figure X, see Table X for statistical outputs). Although also for the action focus, the involved group scored higher, there were no statistical differences (W = 1, p > .05). I found no differences between the action and aim scores for the involved group (W = 1, p > .05) and the control group (W = 1, p > .05).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Looking at the relationship between aim and action focussed self-efficacy beliefs and group membership, I found that the involved group reported significantly higher scores regarding the aim focus, than the control group (W = 8, p < .05, Figure X, see Table X for statistical outputs). Although also for the action focus, the involved group scored higher, there were no statistical differences (W = 1, p > .05). I found no differences between the action and aim scores for the involved group (W = 1, p > .05) and the control group (W = 1, p > .05).",
    2) | Out-Null

Write-Host "done"
